# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) on each
# class/job sheet with the latest pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 186.23529
$ws.Range("I9").Value = 199.54546
$ws.Range("J9").Value = 161.83333
$ws.Range("K9").Value = 199.54546
$ws.Range("L9").Value = 161.83333
$ws.Range("M9").Value = -30.54545999999999
$ws.Range("N9").Value = -499.83333
$ws.Range("H19").Value = 782.2857
$ws.Range("I19").Value = 777.8
$ws.Range("K19").Value = 777.8
$ws.Range("M19").Value = -602.8
$ws.Range("H70").Value = 1999.4
$ws.Range("J70").Value = 1999.25
$ws.Range("L70").Value = 5997.75
$ws.Range("N70").Value = -6537.75
$ws.Range("H73").Value = 1999.4
$ws.Range("J73").Value = 1999.25
$ws.Range("L73").Value = 5997.75
$ws.Range("N73").Value = -7869.75
$ws.Range("H115").Value = 560.4
$ws.Range("I115").Value = 560.4
$ws.Range("K115").Value = 1681.2
$ws.Range("M115").Value = -114.1999999999998
$ws.Range("H132").Value = 7443.136
$ws.Range("I132").Value = 7112.95
$ws.Range("K132").Value = 21338.85
$ws.Range("M132").Value = -18808.85
$ws.Range("H137").Value = 5800.6875
$ws.Range("I137").Value = 1858.3334
$ws.Range("J137").Value = 9279.235000000001
$ws.Range("K137").Value = 5575.0002
$ws.Range("L137").Value = 27837.705
$ws.Range("M137").Value = -3025.0002
$ws.Range("N137").Value = -32937.705
$ws.Range("H138").Value = 5137.4355
$ws.Range("J138").Value = 5847.875
$ws.Range("L138").Value = 17543.625
$ws.Range("N138").Value = -27823.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2104
$ws.Range("I32").Value = 2104
$ws.Range("K32").Value = 2104
$ws.Range("M32").Value = -1817
$ws.Range("H37").Value = 45821.6
$ws.Range("I37").Value = 25034
$ws.Range("J37").Value = 51018.5
$ws.Range("K37").Value = 25034
$ws.Range("L37").Value = 51018.5
$ws.Range("M37").Value = -24761
$ws.Range("N37").Value = -51564.5
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H63").Value = 2249.5
$ws.Range("I63").Value = 2249.5
$ws.Range("K63").Value = 2249.5
$ws.Range("M63").Value = -1563.5
$ws.Range("H66").Value = 2249.5
$ws.Range("I66").Value = 2249.5
$ws.Range("K66").Value = 11247.5
$ws.Range("M66").Value = -7815.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 90909430
$ws.Range("I80").Value = 333333440
$ws.Range("J80").Value = 421
$ws.Range("K80").Value = 333333440
$ws.Range("L80").Value = 421
$ws.Range("M80").Value = -333332442
$ws.Range("N80").Value = -2417
$ws.Range("H82").Value = 8056.6665
$ws.Range("I82").Value = 2828.5
$ws.Range("J82").Value = 49882
$ws.Range("K82").Value = 2828.5
$ws.Range("L82").Value = 49882
$ws.Range("M82").Value = -2445.5
$ws.Range("N82").Value = -50648
$ws.Range("H83").Value = 90909430
$ws.Range("I83").Value = 333333440
$ws.Range("J83").Value = 421
$ws.Range("K83").Value = 1666667200
$ws.Range("L83").Value = 2105
$ws.Range("M83").Value = -1666662208
$ws.Range("N83").Value = -12089
$ws.Range("H85").Value = 8056.6665
$ws.Range("I85").Value = 2828.5
$ws.Range("J85").Value = 49882
$ws.Range("K85").Value = 2828.5
$ws.Range("L85").Value = 49882
$ws.Range("M85").Value = -1502.5
$ws.Range("N85").Value = -52534
$ws.Range("H95").Value = 76996
$ws.Range("J95").Value = 76996
$ws.Range("L95").Value = 76996
$ws.Range("N95").Value = -82488
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H138").Value = 64997.273
$ws.Range("J138").Value = 64997.273
$ws.Range("L138").Value = 64997.273
$ws.Range("N138").Value = -75277.273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2656.9126
$ws.Range("I31").Value = 2163.5833
$ws.Range("J31").Value = 2743.9707
$ws.Range("K31").Value = 2163.5833
$ws.Range("L31").Value = 2743.9707
$ws.Range("M31").Value = -1868.5833
$ws.Range("N31").Value = -3333.9707
$ws.Range("H34").Value = 2656.9126
$ws.Range("I34").Value = 2163.5833
$ws.Range("J34").Value = 2743.9707
$ws.Range("K34").Value = 2163.5833
$ws.Range("L34").Value = 2743.9707
$ws.Range("M34").Value = -1961.5833
$ws.Range("N34").Value = -3147.9707
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H75").Value = 44998
$ws.Range("J75").Value = 44998
$ws.Range("L75").Value = 44998
$ws.Range("N75").Value = -46994
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H78").Value = 44998
$ws.Range("J78").Value = 44998
$ws.Range("L78").Value = 134994
$ws.Range("N78").Value = -144978
$ws.Range("H94").Value = 385.66666
$ws.Range("J94").Value = 107.333336
$ws.Range("L94").Value = 107.333336
$ws.Range("N94").Value = -1009.333336
$ws.Range("H127").Value = 36560.75
$ws.Range("J127").Value = 36560.75
$ws.Range("L127").Value = 36560.75
$ws.Range("N127").Value = -46480.75
$ws.Range("H134").Value = 3524.611
$ws.Range("I134").Value = 3456.4482
$ws.Range("J134").Value = 3807
$ws.Range("K134").Value = 10369.3446
$ws.Range("L134").Value = 11421
$ws.Range("M134").Value = -7834.3446
$ws.Range("N134").Value = -16491

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 1999
$ws.Range("J31").Value = 1999
$ws.Range("L31").Value = 5997
$ws.Range("N31").Value = -6573
$ws.Range("H55").Value = 5646.385
$ws.Range("I55").Value = 2760
$ws.Range("J55").Value = 7450.375
$ws.Range("K55").Value = 8280
$ws.Range("L55").Value = 22351.125
$ws.Range("M55").Value = -8103
$ws.Range("N55").Value = -22705.125
$ws.Range("H68").Value = 1698049.9
$ws.Range("I68").Value = 1356.5333
$ws.Range("J68").Value = 2276468
$ws.Range("K68").Value = 4069.5999
$ws.Range("L68").Value = 6829404
$ws.Range("M68").Value = -3258.5999
$ws.Range("N68").Value = -6831026
$ws.Range("H71").Value = 1698049.9
$ws.Range("I71").Value = 1356.5333
$ws.Range("J71").Value = 2276468
$ws.Range("K71").Value = 12208.7997
$ws.Range("L71").Value = 20488212
$ws.Range("M71").Value = -8152.7997
$ws.Range("N71").Value = -20496324
$ws.Range("H107").Value = 7801.1294
$ws.Range("J107").Value = 8186.569
$ws.Range("L107").Value = 24559.707
$ws.Range("N107").Value = -28399.707

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13422.77
$ws.Range("I70").Value = 4699.5
$ws.Range("K70").Value = 4699.5
$ws.Range("M70").Value = -4429.5
$ws.Range("H73").Value = 13422.77
$ws.Range("I73").Value = 4699.5
$ws.Range("K73").Value = 4699.5
$ws.Range("M73").Value = -3763.5
$ws.Range("H122").Value = 100005230
$ws.Range("I122").Value = 3667.2
$ws.Range("K122").Value = 11001.6
$ws.Range("M122").Value = -8551.599999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2000000
$ws.Range("J2").Value = 2000000
$ws.Range("L2").Value = 2000000
$ws.Range("N2").Value = -2000224
$ws.Range("H22").Value = 357143680
$ws.Range("I22").Value = 35715536
$ws.Range("J22").Value = 1000000000
$ws.Range("K22").Value = 35715536
$ws.Range("L22").Value = 1000000000
$ws.Range("M22").Value = -35715241
$ws.Range("N22").Value = -1000000590
$ws.Range("H27").Value = 357143680
$ws.Range("I27").Value = 35715536
$ws.Range("J27").Value = 1000000000
$ws.Range("K27").Value = 35715536
$ws.Range("L27").Value = 1000000000
$ws.Range("M27").Value = -35715429
$ws.Range("N27").Value = -1000000214
$ws.Range("H40").Value = 36984.945
$ws.Range("I40").Value = 48656.08
$ws.Range("K40").Value = 48656.08
$ws.Range("M40").Value = -48520.08
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H127").Value = 100000
$ws.Range("J127").Value = 100000
$ws.Range("L127").Value = 100000
$ws.Range("N127").Value = -109920
$ws.Range("H132").Value = 11518.64
$ws.Range("I132").Value = 10648.3
$ws.Range("K132").Value = 31944.9
$ws.Range("M132").Value = -29414.9
$ws.Range("H136").Value = 6380.1
$ws.Range("I136").Value = 6360.4
$ws.Range("K136").Value = 19081.2
$ws.Range("M136").Value = -16531.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 13899979
$ws.Range("J136").Value = 19000
$ws.Range("L136").Value = 57000
$ws.Range("N136").Value = -62100
$ws.Range("H137").Value = 124170.25
$ws.Range("J137").Value = 124170.25
$ws.Range("L137").Value = 124170.25
$ws.Range("N137").Value = -134370.25
$ws.Range("H141").Value = 276904.34
$ws.Range("J141").Value = 276904.34
$ws.Range("L141").Value = 276904.34
$ws.Range("N141").Value = -287264.34
